$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "Testcodeword5",
    "Testcodeword6",
    "Testcodeword7",
    "Testcodeword8",
    "Testcodeword9",
    "Testcodeword10",
    "Testcodeword11",
    "Testcodeword12"
)

$row = 5
foreach ($word in $words) {
    $ws.Cells.Item($row, 1).Value = $word
    $row++
}

$ws.Range("G8").Select()
